# Apply updated cryptocurrency price/volume data to the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.890.38"
$ws.Range("E2").Value = "  -1.14%  "

$ws.Range("D3").Value = "1.640.32"
$ws.Range("E3").Value = "  -0.77%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.55%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.49%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5030"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.76%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2575"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.74%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06381"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.96%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.57"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.22%  "

$ws.Range("E11").Value = "  -1.32%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.789.46"
$ws.Range("E12").Value = "  +8.19%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.260"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.83%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "1.863.39"
$ws.Range("E14").Value = "  -1.05%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5456"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.06%  "

$ws.Range("D16").Value = "0.0₅7894"
$ws.Range("E16").Value = "  -1.53%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.07"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.15%  "

$ws.Range("D18").Value = "25.876.56"
$ws.Range("E18").Value = "  -1.39%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.003"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.44%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "202.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.73%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.400"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.11%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.885"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.93%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.973"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.76%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.002"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.68%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.885"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.29%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.69"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.89%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1135"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.30%  "

$ws.Range("E28").Value = "  -1.11%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.791"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.99%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.236"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.53%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04974"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.47%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.270"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.70%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.192"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.29%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.544"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.79%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.363"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.32%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.628"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.95%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.8909"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.08%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5640"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.40%  "

$ws.Range("D39").Value = "1.148.03"
$ws.Range("E39").Value = "  -1.73%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01563"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.66%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.001"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.74%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.680"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.19%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.86"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.41%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8062"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.72%  "

$ws.Range("D45").Value = "1.775.07"
$ws.Range("E45").Value = "  -1.19%  "

$ws.Range("E46").Value = "  +4.01%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4535"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.17%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.005"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.16%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.82"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.97%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05051"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.54%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9995"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.83%  "
